$wb = $excel.ActiveWorkbook

$wsList1  = $wb.Worksheets.Item("List1")
$wsSheet1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Sheet1 ("Sheet1" tab, stored as xl/worksheets/sheet2.xml)
# New tuning-reference data block in columns I/J
# ---------------------------------------------------------------------
$wsSheet1.Range("I2").Value = 440
$wsSheet1.Range("J2").Value = 400000

$wsSheet1.Range("I3").Formula = "=1/I2"
$wsSheet1.Range("J3").Formula = "=1/J2"

$wsSheet1.Range("I6").Formula = "=I3/J3"
$wsSheet1.Range("J6").Formula = "=I6*4"

$wsSheet1.Range("I7").Formula = "=I6*`$C6"
$wsSheet1.Range("J7").Formula = "=J6*`$C6"
$wsSheet1.Range("I7:J7").ClearFormats() | Out-Null

# Selection on the Sheet1 tab is now I7 (the tab is the active one)
$wsSheet1.Range("I7").Select() | Out-Null

# ---------------------------------------------------------------------
# List1 tab (xl/worksheets/sheet1.xml)
# Column P was widened (best fit) and a range was selected there
# ---------------------------------------------------------------------
$wsList1.Columns.Item(16).ColumnWidth = 9.140625

$wsList1.Range("O2:P6").Select() | Out-Null

# Make Sheet1 the active sheet/tab again (it is the last-selected tab)
$wsSheet1.Activate() | Out-Null
$wsSheet1.Range("I7").Select() | Out-Null
